$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2021" data point was added to the table, to the right of the
# existing 2007-2020 series (columns D-Q). Extend column R with the same
# look as column Q (same row) and fill in the new figures.

$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("R4").Value = 2021

$ws.Range("R5").Value = 5.8
$ws.Range("R6").Value = 4.7
$ws.Range("R7").Value = 1.6
$ws.Range("R8").Value = 12.9
$ws.Range("R9").Value = 10.2
$ws.Range("R10").Value = 4.2
$ws.Range("R11").Value = 3.3
$ws.Range("R12").Value = 15.2
$ws.Range("R13").Value = 2.4
$ws.Range("R14").Value = 0.6

[void]$ws.Range("T9").Select()
